$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.324.29'
$ws.Range('E2').Value = '  -1.40%  '

# Row 3
$ws.Range('D3').Value = '1.590.28'
$ws.Range('E3').Value = '  -0.47%  '

# Row 4
$ws.Range('E4').Value = '  -0.58%  '

# Row 5
$ws.Range('D5').Value = "'209.89"
$ws.Range('D5').Style = 'Normal'

# Row 6
$ws.Range('E6').Value = '  -1.39%  '

# Row 7
$ws.Range('E7').Value = '  -0.49%  '

# Row 8
$ws.Range('E8').Value = '  -1.22%  '

# Row 9
$ws.Range('E9').Value = '  -0.69%  '

# Row 10
$ws.Range('D10').Value = "'19.53"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.86%  '

# Row 11
$ws.Range('D11').Value = "'0.0845"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.20%  '

# Row 12
$ws.Range('D12').Value = '1.815.07'
$ws.Range('E12').Value = '  -0.43%  '

# Row 13
$ws.Range('D13').Value = '1.617.12'
$ws.Range('E13').Value = '  +0.68%  '

# Row 14
$ws.Range('E14').Value = '  +0.56%  '

# Row 15
$ws.Range('E15').Value = '  -1.49%  '

# Row 16
$ws.Range('D16').Value = "'64.53"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.65%  '

# Row 17
$ws.Range('D17').Value = '26.340.04'

# Row 18
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  -1.87%  '

# Row 19
$ws.Range('E19').Value = '  +4.82%  '

# Row 20
$ws.Range('D20').Value = "'210.93"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.77%  '

# Row 21
$ws.Range('E21').Value = '  -0.44%  '

# Row 22
$ws.Range('E22').Value = '  -0.65%  '

# Row 23
$ws.Range('D23').Value = "'2.17"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.15%  '

# Row 24
$ws.Range('D24').Value = "'8.91"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.97%  '

# Row 25
$ws.Range('D25').Value = "'145.11"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.38%  '

# Row 26
$ws.Range('E26').Value = '  -0.46%  '

# Row 27
$ws.Range('D27').Value = "'7.05"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.31%  '

# Row 28
$ws.Range('E28').Value = '  -0.74%  '

# Row 29
$ws.Range('E29').Value = '  -0.30%  '

# Row 30
$ws.Range('D30').Value = "'0.0504"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.62%  '

# Row 31
$ws.Range('D31').Value = "'1.14"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.98%  '

# Row 32
$ws.Range('E32').Value = '  -1.52%  '

# Row 33
$ws.Range('E33').Value = '  +0.11%  '

# Row 34
$ws.Range('D34').Value = '1.304.84'
$ws.Range('E34').Value = '  +2.57%  '

# Row 35
$ws.Range('E35').Value = '  +2.90%  '

# Row 36
$ws.Range('E36').Value = '  -1.84%  '

# Row 37
$ws.Range('E37').Value = '  -0.79%  '

# Row 38
$ws.Range('E38').Value = '  -0.58%  '

# Row 39
$ws.Range('D39').Value = "'1.10"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.19%  '

# Row 40
$ws.Range('D40').Value = "'0.812"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.56%  '

# Row 41
$ws.Range('E41').Value = '  -0.38%  '

# Row 42
$ws.Range('E42').Value = '  +3.48%  '

# Row 43
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = "'62.66"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.06%  '

# Row 44
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = "'2.14"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.53%  '

# Row 45
$ws.Range('E45').Value = '  -1.83%  '

# Row 46
$ws.Range('D46').Value = '1.726.92'
$ws.Range('E46').Value = '  -0.43%  '

# Row 47
$ws.Range('D47').Value = "'87.99"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.70%  '

# Row 48
$ws.Range('E48').Value = '  -4.80%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = "'0.0981"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.63%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.0504"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.65%  '

# Row 51
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = "'1.00"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.48%  '
